# 7.10 Fixed Some Bugs
# Clear the leftover/erroneous H4:I4 values (H4 had a stray 0, I4 referenced
# the unused "Dee-Thinking1" string) and leave the selection on that range,
# matching the author's cleanup pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H4:I4").ClearContents()
$ws.Range("H4:I4").Select()
